$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (shifts THIAGO and everything below down by one)
$ws.Rows.Item(2).Insert()

# Fill in the new row 2 with the new record
# (Force the account number to be stored as text so the leading zeros
#  survive, then strip the formatting change back off so the cell keeps
#  the workbook's default style, matching the other data rows.)
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "005366671"
$ws.Cells.Item(2, 1).ClearFormats()
$ws.Cells.Item(2, 2).Value = "TATIANA"
$ws.Cells.Item(2, 3).Value = 31900

# The ANDRE row (originally row 3) has now shifted down to row 4 - remove it
$ws.Rows.Item(4).Delete()
